$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.522586584091187
$ws.Range("B1").Value = 2.050639152526855
$ws.Range("C1").Value = 2.291865348815918
$ws.Range("D1").Value = 2.807581901550293
$ws.Range("E1").Value = 2.608596563339233
